# Update the datetime8 field text ("last printed"/auto-update date field cache)
# from "9/17/18 3:42 PM" to "3/4/19 8:29 PM" across the Handout Master, Notes
# Master, and every slide's Notes Page that carries a Date placeholder.

$p = $ppt.ActivePresentation
$oldText = "9/17/18 3:42 PM"
$newText = "3/4/19 8:29 PM"

function Update-DateField($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldText) {
                $tr.Text = $newText
            }
        }
    }
}

# Handout master date placeholder
Update-DateField $p.HandoutMaster.Shapes

# Notes master date placeholder
Update-DateField $p.NotesMaster.Shapes

# Every slide's notes page date placeholder (only slides with a notes page
# relationship have one; others are skipped automatically since their text
# won't match).
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    Update-DateField $s.NotesPage.Shapes
}
